$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original styles for columns D and E, and force text format
# so that numeric-looking values are not auto-converted by Excel.
$styleD = $ws.Range("D2:D51").Style
$styleE = $ws.Range("E2:E51").Style
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '65.762.29'
$ws.Range("E2").Value = '  -3.20%  '
$ws.Range("D3").Value = '3.514.16'
$ws.Range("E3").Value = '  -0.05%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '557.33'
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").Value = '180.57'
$ws.Range("E6").Value = '  -5.71%  '
$ws.Range("D7").Value = '0.640'
$ws.Range("E7").Value = '  +3.95%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '0.635'
$ws.Range("E9").Value = '  -1.23%  '
$ws.Range("D10").Value = '0.153'
$ws.Range("E10").Value = '  +2.19%  '
$ws.Range("D11").Value = '54.09'
$ws.Range("E11").Value = '  -6.44%  '
$ws.Range("D12").Value = '0.0000273'
$ws.Range("E12").Value = '  -1.87%  '
$ws.Range("D13").Value = '9.30'
$ws.Range("E13").Value = '  -2.49%  '
$ws.Range("D14").Value = '4.062.29'
$ws.Range("E14").Value = '  +0.26%  '
$ws.Range("D15").Value = '3.499.32'
$ws.Range("E15").Value = '  +0.12%  '
$ws.Range("D16").Value = '18.59'
$ws.Range("E16").Value = '  +0.54%  '
$ws.Range("E17").Value = '  +0.26%  '
$ws.Range("D18").Value = '12.17'
$ws.Range("E18").Value = '  +2.30%  '
$ws.Range("D19").Value = '65.734.40'
$ws.Range("E19").Value = '  -3.58%  '
$ws.Range("D20").Value = '0.999'
$ws.Range("E20").Value = '  -1.57%  '
$ws.Range("D21").Value = '419.44'
$ws.Range("E21").Value = '  +3.04%  '
$ws.Range("D22").Value = '4.06'
$ws.Range("E22").Value = '  +2.05%  '
$ws.Range("D23").Value = '86.65'
$ws.Range("E23").Value = '  +1.73%  '
$ws.Range("D24").Value = '4.14'
$ws.Range("E24").Value = '  -1.92%  '
$ws.Range("D25").Value = '12.95'
$ws.Range("E25").Value = '  +8.47%  '
$ws.Range("D26").Value = '10.81'
$ws.Range("E26").Value = '  -11.74%  '
$ws.Range("D27").Value = '2.87'
$ws.Range("E27").Value = '  -2.52%  '
$ws.Range("E28").Value = '  -3.81%  '
$ws.Range("D29").Value = '9.12'
$ws.Range("E29").Value = '  +4.78%  '
$ws.Range("D30").Value = '30.45'
$ws.Range("E30").Value = '  -0.52%  '
$ws.Range("D31").Value = '6.53'
$ws.Range("E31").Value = '  -6.36%  '
$ws.Range("D32").Value = '608.64'
$ws.Range("E32").Value = '  -11.36%  '
$ws.Range("D33").Value = '11.77'
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("E34").Value = '  -0.76%  '
$ws.Range("D35").Value = '59.72'
$ws.Range("E35").Value = '  -2.12%  '
$ws.Range("E36").Value = '  +9.21%  '
$ws.Range("B37").Value = 'InjectiveProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D37").Value = '37.59'
$ws.Range("E37").Value = '  -4.41%  '
$ws.Range("B38").Value = 'Dai'
$ws.Range("C38").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.05%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '3.415.87'
$ws.Range("E39").Value = '  +11.82%  '
$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '0.0₃0794'
$ws.Range("E40").Value = '  -5.45%  '
$ws.Range("D41").Value = '0.382'
$ws.Range("E41").Value = '  -6.08%  '
$ws.Range("E42").Value = '  +0.22%  '
$ws.Range("D43").Value = '3.27'
$ws.Range("E43").Value = '  -4.72%  '
$ws.Range("D44").Value = '2.87'
$ws.Range("E44").Value = '  -3.64%  '
$ws.Range("E45").Value = '  -10.10%  '
$ws.Range("D46").Value = '0.0416'
$ws.Range("E46").Value = '  -2.10%  '
$ws.Range("D47").Value = '3.23'
$ws.Range("E47").Value = '  -0.80%  '
$ws.Range("D48").Value = '2.69'
$ws.Range("E48").Value = '  -2.31%  '
$ws.Range("D49").Value = '0.133'
$ws.Range("E49").Value = '  +1.54%  '
$ws.Range("D50").Value = '8.50'
$ws.Range("E50").Value = '  -4.55%  '
$ws.Range("D51").Value = '138.13'
$ws.Range("E51").Value = '  -1.49%  '

# Restore original styles
$ws.Range("D2:D51").Style = $styleD
$ws.Range("E2:E51").Style = $styleE

Write-Host "Applied crypto list updates"
